# Preparing for Alaska Merge
# - Update the frozen-pane/selection view state on Sheet1
# - Rows 13, 14 and 23 (xdev / Channel / xks) move from the
#   "DAVNET (Guild)" role/req to "CSD", and each gains a Page ("Tools")
#   and Channel ("alaska_classified") value, matching the formatting
#   already used by similar rows (Page->Tools uses the same fill as D4,
#   Channel->alaska_classified uses the same fill as E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet view: frozen pane anchor + active selection -------------------
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("C6").Select() | Out-Null

# --- helper: copy formatting from a template cell onto a target cell -----
function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

$rows = 13, 14, 23

foreach ($r in $rows) {
    Copy-CellFormat "B5" "B$r"
    $ws.Range("B$r").Value = "CSD"

    Copy-CellFormat "D4" "D$r"
    $ws.Range("D$r").Value = "Tools"

    Copy-CellFormat "E7" "E$r"
    $ws.Range("E$r").Value = "alaska_classified"
}

$ws.Application.CutCopyMode = $false

Write-Output "Updated rows 13, 14, 23 and refreshed the sheet view"
